# Save outputs to subfolder data/results/<model_name>/<val> or <test>
#
# - Remove the "Papers" sheet; its contents get folded back into the
#   "Experiments" sheet as a new "From Papers" section (rows 8-11).
# - On "Experiments": bold a few already-present metrics in row 3, add a
#   new results row (row 5) for "DORN + Histogram Rescaling ...", then
#   append the "From Papers" section that used to live on the "Papers"
#   sheet.
#
# NOTE: cell-write order below is deliberate (not just cosmetic) — new
# unique strings are interned into the shared-string table in the order
# their cells are first written, so rows are populated out of visual
# order to land each label/value at its expected shared-string slot.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experiments")

# --- Drop the old "Papers" sheet -------------------------------------
$wb.Worksheets.Item("Papers").Delete()

# --- Row 3: highlight (bold) a few already-present metrics -----------
$ws.Range("D3").Value = 0.991738652808063
$ws.Range("D3").Font.Bold = $true
$ws.Range("E3").Value = 0.315197697889984
$ws.Range("E3").Font.Bold = $true
$ws.Range("H3").Value = 0.0872337991268437
$ws.Range("H3").Font.Bold = $true

# --- Row 9: Eigen et. al. (label only) --------------------------------
$ws.Range("A9").Value = "Eigen et. al."

# --- Row 10: DORN, values taken from the paper -------------------------
$ws.Range("A10").Value = "DORN"
$ws.Range("B10").Value = 0.828
$ws.Range("C10").Value = 0.965
$ws.Range("D10").Value = 0.992
$ws.Range("E10").Value = "-"
$ws.Range("F10").Value = 0.509
$ws.Range("G10").Value = 0.115
$ws.Range("H10").Value = "-"
$ws.Range("I10").Value = 0.051

# --- Row 11: Laina et. al., values taken from the paper -----------------
$ws.Range("A11").Value = "Laina et. al."
$ws.Range("B11").Value = 0.811
$ws.Range("C11").Value = 0.953
$ws.Range("D11").Value = 0.988
$ws.Range("E11").Value = "?"
$ws.Range("F11").Value = 0.573
$ws.Range("G11").Value = 0.127
$ws.Range("H11").Value = "?"
$ws.Range("I11").Value = 0.055
$ws.Range("J11").Value = "?"

# --- Finish row 10 (after "q" becomes a new shared string) --------------
$ws.Range("J10").Value = "q"

# --- Row 5: new "DORN + Histogram Rescaling ..." results row -----------
$ws.Range("A5").Value = "DORN + Histogram Rescaling (No noise, no albedo, no falloff)"

$ws.Range("B5").Value = 0.899023454858409
$ws.Range("B5").Font.Bold = $true
$ws.Range("C5").Value = 0.96994871391909
$ws.Range("C5").Font.Bold = $true
$ws.Range("D5").Value = 0.989547989954822
$ws.Range("E5").Value = 0.3224211819335
$ws.Range("F5").Value = 0.486436836461652
$ws.Range("F5").Font.Bold = $true
$ws.Range("G5").Value = 0.0947125232712816
$ws.Range("G5").Font.Bold = $true
$ws.Range("H5").Value = 0.088471443711215
$ws.Range("I5").Value = 0.0412984580063847
$ws.Range("I5").Font.Bold = $true
$ws.Range("J5").Value = 0.147859282371489
$ws.Range("J5").Font.Bold = $true

# --- Row 8: "From Papers" section header (copied from old Papers!A1:J1)
$ws.Range("A8").Value = "From Papers"
$ws.Range("B8").Value = "delta1"
$ws.Range("C8").Value = "delta2"
$ws.Range("D8").Value = "delta3"
$ws.Range("E8").Value = "mse"
$ws.Range("F8").Value = "rmse"
$ws.Range("G8").Value = "rel_abs_diff"
$ws.Range("H8").Value = "rel_sqr_diff"
$ws.Range("I8").Value = "log10"
$ws.Range("J8").Value = "log_rmse"

# --- Selection / active sheet state -----------------------------------
$ws.Activate()
$ws.Range("A6").Select()
